$d = $word.ActiveDocument

# Merge a range of text (found via whole-document Find) that is currently split
# across several runs sharing identical (empty) formatting into a single run,
# restoring the explicit-but-empty <w:rPr/> element that Word keeps on such runs.
function Merge-EmptyRprRun {
    param($doc, $searchText, $finalText)

    $rng = $doc.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $found = $doc.Range($rng.Start, $rng.End)

    # Step 1: make a genuine textual change (append a marker) - this forces Word to
    # coalesce the matched runs into a single run (formatting is lost in the process).
    $found.Text = $finalText + "X"
    $mergedLen = $finalText.Length + 1
    $merged = $doc.Range($found.Start, $found.Start + $mergedLen)

    # Step 2: change back to the real final text (still a single run).
    $merged.Text = $finalText
    $finalLen = $finalText.Length
    $finalRange = $doc.Range($found.Start, $found.Start + $finalLen)

    # Step 3: toggle a property that is not otherwise used in the document, to force
    # Word to re-materialize an explicit (empty) <w:rPr/> on the run.
    $finalRange.Bold = 1
    $finalRange.Bold = 0
}

# Merge a range of text that is split across several runs sharing identical
# *non-empty* formatting (e.g. bold, or color/size) into one run. A plain
# Find/Replace is enough here because Word preserves non-empty rPr content
# when coalescing runs.
function Merge-FormattedRun {
    param($doc, $searchText, $finalText)
    $doc.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $finalText, 2)
}

# 1. Title: "Professor Vajda Weekly Meeting #" + "2" -> one bold run
Merge-FormattedRun $d "Professor Vajda Weekly Meeting #2" "Professor Vajda Weekly Meeting #2"

# 2. Subtitle: "January " + "17" + "th, 2022 / 10:00 AM / Zoom" -> one run
Merge-FormattedRun $d "January 17th, 2022 / 10:00 AM / Zoom" "January 17th, 2022 / 10:00 AM / Zoom"

# 3-6. Heading3 time ranges, each split across 3 runs with empty rPr -> one run each
Merge-EmptyRprRun $d "10:00 AM - 10:13 AM" "10:00 AM - 10:13 AM"
Merge-EmptyRprRun $d "10:13 AM - 10:18 AM" "10:13 AM - 10:18 AM"
Merge-EmptyRprRun $d "10:18 AM - 10:35 AM" "10:18 AM - 10:35 AM"
Merge-EmptyRprRun $d "10:35 AM - 10:53 AM" "10:35 AM - 10:53 AM"

# 7. "Meet with the client and discuss questions we have to finish our SRS's."
#    becomes 3 runs: "Meet with the client and discuss " / "any " / "questions we have to finish our SRS's."
$rng = $d.Content
$rng.Find.Execute("Meet with the client and discuss questions we have to finish our SRS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found = $d.Range($rng.Start, $rng.End)
$start = $found.Start
$found.Text = "Meet with the client and discuss any questions we have to finish our SRS"

$para = $found.Paragraphs(1)
$paraEnd = $para.Range.End - 1

$r1 = $d.Range($start, $start + 33)
$r2 = $d.Range($start + 33, $start + 37)
$r3 = $d.Range($start + 37, $paraEnd)

$r1.Bold = 1
$r2.Bold = 1
$r3.Bold = 1
$r1.Bold = 0
$r2.Bold = 0
$r3.Bold = 0
